# Regenerate merged AHB files
# 1. Rename header columns: "_old" -> "_FV2310", "_new" -> "_FV2404"
# 2. Convert the used range A1:U73 into an Excel Table ("Table1")
# 3. Freeze panes at row 2 (split below header row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# Create the table over the full data range
$tableRange = $ws.Range("A1:U73")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze panes below the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
